# Apply the changes described by the commit:
#  - Update timestamp embedded in the header text of the "기본설정" (conf) sheet
#  - Clear the value of D3 on the "기본설정" sheet (keeping its style)
#  - Add merged cell ranges A1:B1 and C3:D3 on the "기본설정" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("기본설정")

# Update the header text timestamp (A1 = "해더1-1722325964" -> "해더1-1722391550")
$ws.Range("A1").Value = "해더1-1722391550"

# Clear D3's value while keeping its formatting/style
$ws.Range("D3").ClearContents()

# Merge the required cell ranges
$ws.Range("A1:B1").Merge()
$ws.Range("C3:D3").Merge()
